# edit.ps1 -- Apply the "Update countries & provincias Spain" data refresh
# to paises.xlsx. Only the COVID-19 statistic columns (B:H -- Casos totales,
# Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy,
# Muertes) for a handful of countries changed between the two data pulls,
# plus the "last updated" footer timestamp in A1. Country names/order in
# column A are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 10:28"

# Update statistic cells that changed between the two data pulls
# Row 6 (Rusia)
$ws.Range("B6").Value = 592280
$ws.Range("C6").Value = 7600
$ws.Range("D6").Value = 344416
$ws.Range("E6").Value = 239658
$ws.Range("G6").Value = 95
$ws.Range("H6").Value = 8206
# Row 35 (Singapur)
$ws.Range("B35").Value = 42313
$ws.Range("C35").Value = 218
$ws.Range("E35").Value = 7345
# Row 39 (Polonia)
$ws.Range("D39").Value = 17076
$ws.Range("E39").Value = 13499
# Row 105 (Estonia)
$ws.Range("D105").Value = 1765
$ws.Range("E105").Value = 147
# Row 112 (Lituania)
$ws.Range("B112").Value = 1801
$ws.Range("C112").Value = 3
$ws.Range("E112").Value = 250
# Row 115 (Libano)
$ws.Range("B115").Value = 1588
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 1447
$ws.Range("E115").Value = 113
$ws.Range("H115").Value = 28
# Row 116 (Eslovaquia)
$ws.Range("D116").Value = 1068
$ws.Range("E116").Value = 487
$ws.Range("H116").Value = 32
# Row 130 (Georgia)
$ws.Range("B130").Value = 915
$ws.Range("C130").Value = 82
$ws.Range("D130").Value = 439
$ws.Range("E130").Value = 473
$ws.Range("H130").Value = 3
# Row 131 (Burkina Faso)
$ws.Range("B131").Value = 908
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 761
$ws.Range("E131").Value = 133
$ws.Range("H131").Value = 14
# Row 132 (Cabo Verde)
$ws.Range("B132").Value = 903
$ws.Range("D132").Value = 814
$ws.Range("E132").Value = 36
$ws.Range("H132").Value = 53
# Row 133 (Congo)
$ws.Range("B133").Value = 890
$ws.Range("D133").Value = 413
$ws.Range("E133").Value = 469
$ws.Range("H133").Value = 8
# Row 134 (Uruguay)
$ws.Range("B134").Value = 883
$ws.Range("D134").Value = 391
$ws.Range("E134").Value = 465
$ws.Range("H134").Value = 27
# Row 135 (Republica del Chad)
$ws.Range("B135").Value = 876
$ws.Range("D135").Value = 814
$ws.Range("E135").Value = 37
$ws.Range("H135").Value = 25
# Row 136 (Principado de Andorra)
$ws.Range("B136").Value = 858
$ws.Range("D136").Value = 752
$ws.Range("E136").Value = 32
$ws.Range("H136").Value = 74
# Row 137 (Estado de Palestina)
$ws.Range("B137").Value = 855
$ws.Range("D137").Value = 792
$ws.Range("E137").Value = 11
$ws.Range("H137").Value = 52
# Row 186 (Polinesia Francesa)
$ws.Range("B186").Value = 63
$ws.Range("C186").Value = 8
$ws.Range("D186").Value = 21
$ws.Range("E186").Value = 42
# Row 187 (Namibia)
$ws.Range("B187").Value = 60
$ws.Range("D187").Value = 60
$ws.Range("E187").Value = 0
# Row 209 (Santa Sede)
$ws.Range("C209").Value = 8
$ws.Range("D209").Value = 2
$ws.Range("E209").Value = 10
# Row 210 (Montserrat)
$ws.Range("B210").Value = 12
$ws.Range("D210").Value = 12
$ws.Range("H210").Value = 0
# Row 211 (Seychelles)
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
# Row 212 (Sahara Occidental)
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 11
$ws.Range("E212").Value = 0
$ws.Range("H212").Value = 0
# Row 213 (Papua Nueva Guinea)
$ws.Range("B213").Value = 10
$ws.Range("C213").Value = 1
$ws.Range("E213").Value = 1
$ws.Range("H213").Value = 1
# Row 215 (Bonaire, San Eustaquio y Saba)
$ws.Range("B215").Value = 8
$ws.Range("D215").Value = 8
# Row 216 (San Bartolome)
$ws.Range("B216").Value = 7
$ws.Range("D216").Value = 7
# Row 217 (Lesoto)
$ws.Range("B217").Value = 6
$ws.Range("D217").Value = 6
$ws.Range("E217").Value = 0
